# Fruta / hortaliza, semanal
# Insert a new weekly price-report record for "Achicoria" (Vega Modelo de
# Temuco) as row 26, pushing the existing rows 26-51 down to 27-52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 26 (Excel shifts rows 26:51 down to 27:52 and
# extends the used range to row 52).
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = 44763
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = 100112010
$ws.Cells.Item(26, 7).Value = "Achicoria"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 65
$ws.Cells.Item(26, 11).Value = 11000
$ws.Cells.Item(26, 12).Value = 11000
$ws.Cells.Item(26, 13).Value = 11000
$ws.Cells.Item(26, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(26, 15).Value = "Región Metropolitana"
$ws.Cells.Item(26, 16).Value = 611
$ws.Cells.Item(26, 17).Value = 18
$ws.Cells.Item(26, 18).Value = "Hortaliza"
